# Lattice multiplication exercises: refresh each table cell's worked-problem
# text (header "A x B" line, the two split-digit headers, the divider, and the
# two lattice-row leading digits) to match the newly generated set of problems.
$d = $word.ActiveDocument
$t = $d.Tables(1)
$vt = [char]11   # <w:br/> manual line break, as used between each w:t run in a cell

# Row 1, Col 1: "69 x 60" -> "72 x 76"
$t.Cell(1,1).Range.Text = "72 x 76" + $vt + "  7    6" + $vt + "  ----" + $vt + "7|    |" + $vt + "2|    |"

# Row 1, Col 2: "84 x 43" -> "23 x 43"
$t.Cell(1,2).Range.Text = "23 x 43" + $vt + "  4    3" + $vt + "  ----" + $vt + "2|    |" + $vt + "3|    |"

# Row 1, Col 3: "23 x 97" -> "93 x 40"
$t.Cell(1,3).Range.Text = "93 x 40" + $vt + "  4    0" + $vt + "  ----" + $vt + "9|    |" + $vt + "3|    |"

# Row 2, Col 1: "84 x 22" -> "27 x 75"
$t.Cell(2,1).Range.Text = "27 x 75" + $vt + "  7    5" + $vt + "  ----" + $vt + "2|    |" + $vt + "7|    |"

# Row 2, Col 2: "32 x 59" -> "46 x 22"
$t.Cell(2,2).Range.Text = "46 x 22" + $vt + "  2    2" + $vt + "  ----" + $vt + "4|    |" + $vt + "6|    |"

# Row 2, Col 3: "20 x 16" -> "18 x 92"
$t.Cell(2,3).Range.Text = "18 x 92" + $vt + "  9    2" + $vt + "  ----" + $vt + "1|    |" + $vt + "8|    |"

# Row 3, Col 1: "25 x 33" -> "12 x 64"
$t.Cell(3,1).Range.Text = "12 x 64" + $vt + "  6    4" + $vt + "  ----" + $vt + "1|    |" + $vt + "2|    |"

# Row 3, Col 2: "75 x 18" -> "75 x 39"
$t.Cell(3,2).Range.Text = "75 x 39" + $vt + "  3    9" + $vt + "  ----" + $vt + "7|    |" + $vt + "5|    |"

# Row 3, Col 3: "96 x 85" -> "91 x 91"
$t.Cell(3,3).Range.Text = "91 x 91" + $vt + "  9    1" + $vt + "  ----" + $vt + "9|    |" + $vt + "1|    |"

# Row 4, Col 1: "48 x 44" -> "20 x 79"
$t.Cell(4,1).Range.Text = "20 x 79" + $vt + "  7    9" + $vt + "  ----" + $vt + "2|    |" + $vt + "0|    |"

# Row 4, Col 2: "83 x 26" -> "95 x 38"
$t.Cell(4,2).Range.Text = "95 x 38" + $vt + "  3    8" + $vt + "  ----" + $vt + "9|    |" + $vt + "5|    |"

# Row 4, Col 3: "44 x 37" -> "11 x 12"
$t.Cell(4,3).Range.Text = "11 x 12" + $vt + "  1    2" + $vt + "  ----" + $vt + "1|    |" + $vt + "1|    |"

# Row 5, Col 1: "22 x 29" -> "32 x 66"
$t.Cell(5,1).Range.Text = "32 x 66" + $vt + "  6    6" + $vt + "  ----" + $vt + "3|    |" + $vt + "2|    |"

# Row 5, Col 2: "64 x 30" -> "38 x 26"
$t.Cell(5,2).Range.Text = "38 x 26" + $vt + "  2    6" + $vt + "  ----" + $vt + "3|    |" + $vt + "8|    |"

# Row 5, Col 3: "61 x 77" -> "74 x 84"
$t.Cell(5,3).Range.Text = "74 x 84" + $vt + "  8    4" + $vt + "  ----" + $vt + "7|    |" + $vt + "4|    |"
